$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 12.588
$ws.Range("D3").Value = -7.726999999999999
$ws.Range("A4").Value = -22.09
$ws.Range("C4").Value = -12.927
$ws.Range("D4").Value = -7.727000000000001
$ws.Range("C5").Value = -12.927
$ws.Range("A6").Value = -20.97
$ws.Range("A7").Value = -21.038
$ws.Range("C8").Value = -12.883
$ws.Range("D9").Value = -8.048
$ws.Range("D11").Value = -8.103999999999999
$ws.Range("D14").Value = -8.028
$ws.Range("A16").Value = -21.328
$ws.Range("C16").Value = -12.719
$ws.Range("D18").Value = -7.637
$ws.Range("A20").Value = -22.269
$ws.Range("E20").Value = 12.85
$ws.Range("C22").Value = -12.78
$ws.Range("D25").Value = -7.995
